{"js": "// The Title and Author paragraphs each had their text split across several\n// runs (one per word/space). Collapse each paragraph back down to a single\n// run holding the full text, leaving every other paragraph untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nlet titleParagraph = null;\nlet authorParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (titleParagraph === null && paragraph.style === \"Title\") {\n    titleParagraph = paragraph;\n  } else if (authorParagraph === null && paragraph.style === \"Author\") {\n    authorParagraph = paragraph;\n  }\n  if (titleParagraph !== null && authorParagraph !== null) {\n    break;\n  }\n}\n\n// Fall back to positional lookup in case the styles above aren't present.\nif (titleParagraph === null) {\n  titleParagraph = paragraphs.items[0];\n}\nif (authorParagraph === null) {\n  authorParagraph = paragraphs.items[1];\n}\n\ntitleParagraph.insertText(\"Sigma Notation: Answers\", Word.InsertLocation.replace);\nauthorParagraph.insertText(\"Ifan Howells-Baines, Mark Toner\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The Title and Author paragraphs each had their text split across several\n# runs (one per word/space). Collapse each paragraph back down to a single\n# run holding the full text, leaving every other paragraph untouched.\n$d = $word.ActiveDocument\n\n$titlePara = $null\n$authorPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($titlePara -eq $null -and $p.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n    } elseif ($authorPara -eq $null -and $p.Style.NameLocal -eq \"Author\") {\n        $authorPara = $p\n    }\n    if ($titlePara -ne $null -and $authorPara -ne $null) {\n        break\n    }\n}\n\n# Fall back to positional lookup in case the styles above aren't present.\nif ($titlePara -eq $null) {\n    $titlePara = $d.Paragraphs.Item(1)\n}\nif ($authorPara -eq $null) {\n    $authorPara = $d.Paragraphs.Item(2)\n}\n\n$titleText = \"Sigma Notation: Answers\"\n$titleRange = $titlePara.Range\n$titleSearch = $titleRange.Text.TrimEnd([char]13, [char]7)\n$titleRange.Find.Execute($titleSearch, $false, $false, $false, $false, $false, $true, 1, $false, $titleText, 2)\n\n$authorText = \"Ifan Howells-Baines, Mark Toner\"\n$authorRange = $authorPara.Range\n$authorSearch = $authorRange.Text.TrimEnd([char]13, [char]7)\n$authorRange.Find.Execute($authorSearch, $false, $false, $false, $false, $false, $true, 1, $false, $authorText, 2)\n"}
